$wb = $excel.ActiveWorkbook

# --- Sheet 1: "ATS Accuracy" ---
$ws1 = $wb.Worksheets.Item("ATS Accuracy")

$ws1.Range("B2").Value = 3
$ws1.Range("C2").Value = 79
$ws1.Range("D2").Value = 82
$ws1.Range("E2").Value = 96.3

$ws1.Range("B3").Value = 2
$ws1.Range("D3").Value = 61
$ws1.Range("E3").Value = 96.7

$ws1.Range("B5").Value = 3
$ws1.Range("C5").Value = 5
$ws1.Range("D5").Value = 8
$ws1.Range("E5").Value = 62.5

$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 7
$ws1.Range("E6").Value = 42.9

# --- Sheet 2: "Total Accuracy" ---
$ws2 = $wb.Worksheets.Item("Total Accuracy")

$ws2.Range("B2").Value = 4
$ws2.Range("C2").Value = 67
$ws2.Range("D2").Value = 71
$ws2.Range("E2").Value = 94.40000000000001

$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 64
$ws2.Range("D3").Value = 66
$ws2.Range("E3").Value = 97

$ws2.Range("B4").Value = 3
$ws2.Range("C4").Value = 20
$ws2.Range("D4").Value = 23
$ws2.Range("E4").Value = 87

$ws2.Range("B5").Value = 3
$ws2.Range("D5").Value = 11
$ws2.Range("E5").Value = 72.7

$ws2.Range("B6").Value = 1
$ws2.Range("D6").Value = 3
$ws2.Range("E6").Value = 66.7
